$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the "Sandia" weekly block (before existing row 555),
# shifting the existing rows 555:590 down to 558:593.
$ws.Rows("555:557").Insert()

# Fill in the 3 newly inserted rows with this week's data (same shape as the
# surrounding rows: Mercado ID, Mercado, Region, Fecha, Codreg, Categoria ID,
# Categoria, Variedad, Calidad, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Unidad de comercializacion, Origen, Precio $/Kg,
# Kg o Unidades, Clasificacion).

# Row 555: Extra
$ws.Cells.Item(555, 1).Value = 3
$ws.Cells.Item(555, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(555, 3).Value = "Coquimbo"
$ws.Cells.Item(555, 4).Value = 44931
$ws.Cells.Item(555, 5).Value = 5
$ws.Cells.Item(555, 6).Value = 100112028
$ws.Cells.Item(555, 7).Value = "Sandia"
$ws.Cells.Item(555, 8).Value = "Sin especificar"
$ws.Cells.Item(555, 9).Value = "Extra"
$ws.Cells.Item(555, 10).Value = 350
$ws.Cells.Item(555, 11).Value = 4000
$ws.Cells.Item(555, 12).Value = 4000
$ws.Cells.Item(555, 13).Value = 4000
$ws.Cells.Item(555, 14).Value = "$/unidad"
$ws.Cells.Item(555, 15).Value = "Paine"
$ws.Cells.Item(555, 16).Value = 4000
$ws.Cells.Item(555, 17).Value = 1
$ws.Cells.Item(555, 18).Value = "Hortaliza"

# Row 556: Primera
$ws.Cells.Item(556, 1).Value = 3
$ws.Cells.Item(556, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(556, 3).Value = "Coquimbo"
$ws.Cells.Item(556, 4).Value = 44931
$ws.Cells.Item(556, 5).Value = 5
$ws.Cells.Item(556, 6).Value = 100112028
$ws.Cells.Item(556, 7).Value = "Sandia"
$ws.Cells.Item(556, 8).Value = "Sin especificar"
$ws.Cells.Item(556, 9).Value = "Primera"
$ws.Cells.Item(556, 10).Value = 400
$ws.Cells.Item(556, 11).Value = 3000
$ws.Cells.Item(556, 12).Value = 3000
$ws.Cells.Item(556, 13).Value = 3000
$ws.Cells.Item(556, 14).Value = "$/unidad"
$ws.Cells.Item(556, 15).Value = "Paine"
$ws.Cells.Item(556, 16).Value = 3000
$ws.Cells.Item(556, 17).Value = 1
$ws.Cells.Item(556, 18).Value = "Hortaliza"

# Row 557: Segunda
$ws.Cells.Item(557, 1).Value = 3
$ws.Cells.Item(557, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(557, 3).Value = "Coquimbo"
$ws.Cells.Item(557, 4).Value = 44931
$ws.Cells.Item(557, 5).Value = 5
$ws.Cells.Item(557, 6).Value = 100112028
$ws.Cells.Item(557, 7).Value = "Sandia"
$ws.Cells.Item(557, 8).Value = "Sin especificar"
$ws.Cells.Item(557, 9).Value = "Segunda"
$ws.Cells.Item(557, 10).Value = 410
$ws.Cells.Item(557, 11).Value = 2000
$ws.Cells.Item(557, 12).Value = 2000
$ws.Cells.Item(557, 13).Value = 2000
$ws.Cells.Item(557, 14).Value = "$/unidad"
$ws.Cells.Item(557, 15).Value = "Paine"
$ws.Cells.Item(557, 16).Value = 2000
$ws.Cells.Item(557, 17).Value = 1
$ws.Cells.Item(557, 18).Value = "Hortaliza"
